$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (28) down to the new row (29)
$ws.Range("A28:H28").Copy()
$ws.Range("A29:H29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row's values
$ws.Cells.Item(29, 1).Value = "2025-08-18 13:09:25 UTC"
$ws.Cells.Item(29, 2).Value = "2025-08-18 18:39:25 IST"
$ws.Cells.Item(29, 3).Value = "SKIPPED"
$ws.Cells.Item(29, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item(29, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item(29, 6).Value = ""
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(29, 8).Value = ""
